# Add data for the "Reverse" for-loop correlated mixed error dataset,
# plus its associated figures/summary values, to the second data block
# (rows 29-38) of Sheet1, in columns Q:W.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New section header in Q29 identifying this as the "Reverse" data set.
$ws.Range("Q29").Value = "Reverse"

# Row 30 - combCor05-04-10
$ws.Range("Q30").Value = "combCor05-04-10"
$ws.Range("R30").Formula = "=0.001241"
$ws.Range("S30").Formula = "=0.00024"
$ws.Range("T30").Formula = "=0.002804"
$ws.Range("U30").Formula = "=0.000382"
$ws.Range("V30").Value = 25.4
$ws.Range("W30").Value = 26.2

# Row 31 - combCor20-10-10
$ws.Range("Q31").Value = "combCor20-10-10"
$ws.Range("R31").Formula = "=0.004597"
$ws.Range("S31").Formula = "=0.000591"
$ws.Range("T31").Formula = "=0.0056925"
$ws.Range("U31").Formula = "=0.000937"
$ws.Range("V31").Value = 25.4
$ws.Range("W31").Value = 26.2

# Row 32 - combCor80-40-10
$ws.Range("Q32").Value = "combCor80-40-10"
$ws.Range("R32").Formula = "=0.020387"
$ws.Range("S32").Formula = "=0.003345"
$ws.Range("T32").Formula = "=0.0202259"
$ws.Range("U32").Formula = "=0.003955"
$ws.Range("V32").Value = 25.4
$ws.Range("W32").Value = 26.2

# Row 33 - combCor05-04-50
$ws.Range("Q33").Value = "combCor05-04-50"
$ws.Range("R33").Formula = "=0.00729"
$ws.Range("S33").Formula = "=0.000763"
$ws.Range("T33").Formula = "=0.004452"
$ws.Range("U33").Formula = "=0.000855"
$ws.Range("V33").Value = 25.4
$ws.Range("W33").Value = 26.2

# Row 34 - combCor20-10-50
$ws.Range("Q34").Value = "combCor20-10-50"
$ws.Range("R34").Formula = "=0.018032"
$ws.Range("S34").Formula = "=0.002525"
$ws.Range("T34").Formula = "=0.016333"
$ws.Range("U34").Formula = "=0.002822"
$ws.Range("V34").Value = 25.4
$ws.Range("W34").Value = 26.2

# Row 35 - combCor80-40-50
$ws.Range("Q35").Value = "combCor80-40-50"
$ws.Range("R35").Formula = "=0.07073"
$ws.Range("S35").Formula = "=0.0048049"
$ws.Range("T35").Formula = "=0.0814612"
$ws.Range("U35").Formula = "=0.01626"
$ws.Range("V35").Value = 25.4
$ws.Range("W35").Value = 26.2

# Row 36 - combCor05-04-80
$ws.Range("Q36").Value = "combCor05-04-80"
$ws.Range("R36").Formula = "=0.005183"
$ws.Range("S36").Formula = "=0.00066"
$ws.Range("T36").Formula = "=0.006492"
$ws.Range("U36").Formula = "=0.001321"
$ws.Range("V36").Value = 25.4
$ws.Range("W36").Value = 26.2

# Row 37 - combCor20-10-80
$ws.Range("Q37").Value = "combCor20-10-80"
$ws.Range("R37").Formula = "=0.0294189"
$ws.Range("S37").Formula = "=0.003344"
$ws.Range("T37").Formula = "=0.019294"
$ws.Range("U37").Formula = "=0.004577"
$ws.Range("V37").Value = 25.4
$ws.Range("W37").Value = 26.2

# Row 38 - combCor80-40-80
$ws.Range("Q38").Value = "combCor80-40-80"
$ws.Range("R38").Formula = "=0.067295"
$ws.Range("S38").Formula = "=0.011521"
$ws.Range("T38").Formula = "=0.118434"
$ws.Range("U38").Formula = "=0.024033"
$ws.Range("V38").Value = 25.7
$ws.Range("W38").Formula = "=12.4+14.2"

# Apply the number formats matching the rest of the sheet: the
# correlated-error columns (R:U) use the 5-decimal custom format,
# while the figure columns (V:W) use the standard 2-decimal format.
$ws.Range("R30:U38").NumberFormat = "0.00000"
$ws.Range("V30:W38").NumberFormat = "0.00"

# Reflect the cursor's final resting place after entering the data.
$ws.Range("R30").Select()
